$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.967.93"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.819.13"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "356.81"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "111.63"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.63%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.559"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +0.00%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.634"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +8.44%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "40.43"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("E11").Value = "  -0.29%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0842"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "20.04"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "3.261.10"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "2.813.49"
$ws.Range("E16").Value = "  +0.94%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.947"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "51.952.97"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("E20").Value = "  +3.48%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.70"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +1.63%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.53"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "268.93"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.77"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "26.26"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  +0.10%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.163"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.06%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "10.46"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.97%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "38.39"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +10.94%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.25"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "6.20"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "52.64"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.68"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +11.14%  "
$ws.Range("E35").Value = "  -0.85%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0875"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.52%  "
$ws.Range("E37").Value = "  +0.03%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "18.92"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +3.38%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "120.51"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "21.97"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.31%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.41"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.111.01"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.42"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +5.83%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.938"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  +9.62%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "5.47"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
